# Auto update Excel log
# Appends new sensor rows to the ALERTS sheet and the mmWave sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ALERTS sheet: add one new CRITICAL / FALL_DETECTED row at row 3
# ---------------------------------------------------------------------
$alerts = $wb.Worksheets.Item("ALERTS")

$alertRow = 3
$alerts.Cells.Item($alertRow, 1).NumberFormat = "@"
$alerts.Cells.Item($alertRow, 1).Value = "2026-01-31"
$alerts.Cells.Item($alertRow, 2).Value = "21:40:49"
$alerts.Cells.Item($alertRow, 3).Value = "21:00"
$alerts.Cells.Item($alertRow, 4).Value = "Living Room"
$alerts.Cells.Item($alertRow, 5).Value = "CRITICAL"
$alerts.Cells.Item($alertRow, 6).Value = "FALL_DETECTED"

# ---------------------------------------------------------------------
# mmWave sheet: add thirteen new PRESENCE_DETECTED / Active rows
# (rows 14-26)
# ---------------------------------------------------------------------
$mmwave = $wb.Worksheets.Item("mmWave")

$timestamps = @(
    "21:40:45",
    "21:40:45",
    "21:40:45",
    "21:40:45",
    "21:40:45",
    "21:40:45",
    "21:40:45",
    "21:40:45",
    "21:40:54",
    "21:41:04",
    "21:41:15",
    "21:41:25",
    "21:41:36"
)

$startRow = 14
for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $r = $startRow + $i
    $mmwave.Cells.Item($r, 1).NumberFormat = "@"
    $mmwave.Cells.Item($r, 1).Value = "2026-01-31"
    $mmwave.Cells.Item($r, 2).Value = $timestamps[$i]
    $mmwave.Cells.Item($r, 3).Value = "21:00"
    $mmwave.Cells.Item($r, 4).Value = "Living Room"
    $mmwave.Cells.Item($r, 5).Value = "PRESENCE_DETECTED"
    $mmwave.Cells.Item($r, 6).Value = "Active"
}
